# Update data-message.xlsx per the commit "Add files via upload":
#  - The "1. Покаяние / 2. Новая Молитва" prayer text gets the section
#    titles wrapped in '#' markers.
#  - The short "Отче наш..." message gets wrapped in '*' markers.
#  - The 21:00-21:04 / 21:05-21:09 triggers attached to the last occurrence
#    of the short message move to 21:55-21:59 / 22:00-22:04.
#  - The sheet view scrolls back to the top and the selection moves to C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prayerText = @"
#1. Покаяние#^
^
Отче наш, Отец Небесный,^
прошу Тебя простить^
все мои прегрешения,^
вольные или невольные.^
^
Аминь!^
^
#2. Новая Молитва#^
^
ОТЧЕ НАШ, ОТЕЦ НЕБЕСНЫЙ,^
Я ПРИНИМАЮ ВЕРУ ТВОЮ,^
ОНА ЕСТЬ МОЙ ПУТЬ.^
^
Я ПРИНИМАЮ КАНОНЫ ТВОИ ВЕЧНЫЕ,^
С ЛЮБОВЬЮ К ТЕБЕ И ДЕЛАМ ТВОИМ,^
ПОДТВЕРЖДАЯ СВОЕЙ ЖИЗНЬЮ^
ВЕРНОСТЬ ТЕБЕ.^
^
ГОСПОДИ, ПРОШУ ДАТЬ МНЕ НАДЕЖДУ^
НА СПАСЕНИЕ ДУШИ МОЕЙ,^
И ДАРОВАТЬ МУДРОСТЬ ТВОЮ^
ДЛЯ ЖИЗНИ МОЕЙ ЗДЕСЬ,^
НА ПЛАНЕТЕ СВЯТАЯ РУСЬ И В ВЕЧНОСТИ.^
^
ПУСТЬ СВЯТА БУДЕТ УВЕРЕННОСТЬ МОЯ,^
ЧТО ТЫ ЕСМЬ!^
^
Господи, я Люблю Тебя, Благодарю Тебя и Уповаю на Милость Твою! Аминь!
"@

$shortText = "*Отче наш, Отец Небесный! Волею Создателя, Пророка и Народа Пространство Святая Русь ЕСМЬ Равенство и Любовь Навечно! Да будет Свет Истины!*"

# Rows 2/4/6/8 hold the long prayer text in column A.
$ws.Range("A2").Value = $prayerText
$ws.Range("A4").Value = $prayerText
$ws.Range("A6").Value = $prayerText
$ws.Range("A8").Value = $prayerText

# Rows 3/5/7/9 hold the short message in column A.
$ws.Range("A3").Value = $shortText
$ws.Range("A5").Value = $shortText
$ws.Range("A7").Value = $shortText
$ws.Range("A9").Value = $shortText

# Replacing the text re-triggers autofit on row 4, whose height was
# manually tuned smaller than the wrapped text needs; restore it so the
# layout stays unchanged. Rows 2/6/8 were already pinned at Excel's row
# height cap (409.5) both before and after, so they need no correction.
$ws.Rows.Item(4).RowHeight = 79.5

# The trigger window tied to the short message's last occurrence shifts
# from 21:00-21:09 to 21:55-22:04.
$ws.Range("B8").Value = "21:55 - 21:59"
$ws.Range("B9").Value = "22:00 - 22:04"

# Reset the view: scroll to the top-left and select C2.
[void]$ws.Range("C2").Select()
